# Update parameter files with WEP scaling and WFF_or_Ben
#
# Adds two new parameter rows to the "Parameters" sheet:
#   Row 56: MFTC_WEP_scaling  = 1   (How should the Winter Energy Payment be scaled? ...)
#   Row 57: WFF_or_Benefit    = Max (What work decision should we assume? ...)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New parameter: MFTC_WEP_scaling -------------------------------------
$ws.Range("B56").Value2 = "MFTC_WEP_scaling"
$ws.Range("C56").Value2 = 1
$ws.Range("D56").Value2 = "How should the Winter Energy Payment be scaled? Average week = 1, Winter week = 12/5, Summer week = 0"

# --- New parameter: WFF_or_Benefit ----------------------------------------
$ws.Range("B57").Value2 = "WFF_or_Benefit"
$ws.Range("C57").Value2 = "Max"
$ws.Range("D57").Value2 = "What work decision should we assume? Go off-benefit and receive IWTC = ""WFF"", stay on-benefit = ""Benefit"", or whichever gives a higher net income = ""Max"""

# --- Formatting -------------------------------------------------------------
# Rows 56:57 re-use the look of the last existing row's Description cell
# (font/fill/left-alignment) but without the bottom border that row 55 has,
# matching the new style added to the workbook (font4/fill3/no border/left).
$ws.Range("D55").Copy()
$newRows = $ws.Range("B56:D57")
$newRows.PasteSpecial(-4122)
$newRows.Borders.LineStyle = -4142

$ws.Range("A1").Select()
